$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.226.94'
$ws.Range("E2").Value = '  -1.17%  '
$ws.Range("D3").Value = '3.335.50'
$ws.Range("E3").Value = '  +2.42%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.70'
$ws.Range("E5").Value = '  -0.78%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '183.69'
$ws.Range("E6").Value = '  +0.45%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  +0.59%  '
$ws.Range("E9").Value = '  -0.24%  '
$ws.Range("E10").Value = '  +0.90%  '
$ws.Range("E11").Value = '  +0.02%  '
$ws.Range("D12").Value = '3.919.37'
$ws.Range("E12").Value = '  +2.42%  '
$ws.Range("E13").Value = '  -0.72%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.27'
$ws.Range("E14").Value = '  +0.29%  '
$ws.Range("D15").Value = '67.440.06'
$ws.Range("E15").Value = '  -0.88%  '
$ws.Range("E16").Value = '  +0.02%  '
$ws.Range("D17").Value = '3.350.12'
$ws.Range("E17").Value = '  +1.77%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '444.27'
$ws.Range("E18").Value = '  +6.96%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.57'
$ws.Range("E19").Value = '  +2.60%  '
$ws.Range("E20").Value = '  -0.81%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.71'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '74.00'
$ws.Range("E22").Value = '  +4.27%  '
$ws.Range("E23").Value = '  -0.08%  '
$ws.Range("D24").Value = '3.490.59'
$ws.Range("E24").Value = '  +2.52%  '
$ws.Range("E25").Value = '  +0.98%  '
$ws.Range("E26").Value = '  +3.35%  '
$ws.Range("E27").Value = '  +3.35%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.02'
$ws.Range("E28").Value = '  -2.47%  '
$ws.Range("E29").Value = '  +0.12%  '
$ws.Range("E30").Value = '  +1.76%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '22.88'
$ws.Range("E31").Value = '  +1.52%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.31'
$ws.Range("E32").Value = '  -1.94%  '
$ws.Range("E33").Value = '  -0.04%  '
$ws.Range("E34").Value = '  -0.12%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.76'
$ws.Range("E35").Value = '  -0.87%  '
$ws.Range("B36").Value = 'Monero'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '161.59'
$ws.Range("E36").Value = '  -1.62%  '
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.49'
$ws.Range("E37").Value = '  +4.34%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '27.48'
$ws.Range("E38").Value = '  +4.10%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.84'
$ws.Range("E39").Value = '  -1.48%  '
$ws.Range("D40").Value = '2.831.55'
$ws.Range("E40").Value = '  +8.13%  '
$ws.Range("E41").Value = '  +0.25%  '
$ws.Range("E42").Value = '  +1.14%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.23'
$ws.Range("E43").Value = '  -0.51%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '40.33'
$ws.Range("E44").Value = '  -0.55%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0670'
$ws.Range("E45").Value = '  -0.01%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '24.44'
$ws.Range("E46").Value = '  +1.56%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.34'
$ws.Range("E47").Value = '  -2.09%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '321.62'
$ws.Range("E48").Value = '  -3.42%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0273'
$ws.Range("E49").Value = '  +0.54%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.984'
$ws.Range("E50").Value = '  +0.87%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '30.91'
$ws.Range("E51").Value = '  +1.84%  '
